# Updates cryptos list figures (prices / 1h volume change) as per the
# "Updated cryptos list" GitHub Actions commit.
#
# NOTE: several "Price" values are plain decimal numbers (e.g. "98.89").
# Excel's Range.Value setter auto-detects such strings as numbers, which
# would change the underlying cell type from text to numeric and break
# the exact text match expected by the workbook (prices are stored as
# text, e.g. so values like "43.072.50" - which aren't valid numbers -
# stay consistent with plain numeric-looking ones like "98.89"). We
# force those particular assignments to remain text by prefixing the
# value with a leading apostrophe, Excel's classic "store as text" cue.
# The apostrophe itself is not stored as part of the cell's text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.072.50"
$ws.Range("E2").Value = "  +0.01%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.312.28"
$ws.Range("E3").Value = "  +0.12%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  -0.10%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'98.89"
$ws.Range("E6").Value = "  -3.04%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.520"
$ws.Range("E7").Value = "  +3.04%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.39%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'35.75"
$ws.Range("E10").Value = "  +0.15%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.77%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.61%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "'17.90"
$ws.Range("E13").Value = "  -0.35%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'6.94"
$ws.Range("E14").Value = "  +0.24%  "

# Row 15 - Wrapped liquid staked Ether 2.0
$ws.Range("D15").Value = "2.672.25"
$ws.Range("E15").Value = "  -0.52%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.253.23"
$ws.Range("E16").Value = "  -2.25%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -3.00%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "42.984.34"
$ws.Range("E18").Value = "  -0.07%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").Value = "'13.60"
$ws.Range("E19").Value = "  +7.20%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  +0.70%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.01%  "

# Row 22 - Litecoin
$ws.Range("E22").Value = "  +0.33%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'240.02"
$ws.Range("E23").Value = "  +1.12%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -2.40%  "

# Row 25 & 26 swap - Dai <-> PancakeSwap (row 25 now holds PancakeSwap, row 26 now holds Dai)
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +0.44%  "

# Row 28 - Monero
$ws.Range("D28").Value = "'168.28"
$ws.Range("E28").Value = "  +0.29%  "

# Row 29 & 30 swap - Toncoin <-> Cosmos (row 29 now holds Cosmos, row 30 now holds Toncoin)
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'9.18"
$ws.Range("E29").Value = "  -0.87%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.05"
$ws.Range("E30").Value = "  -10.77%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").Value = "'33.42"
$ws.Range("E31").Value = "  -3.56%  "

# Row 32 & 33 swap - RenderToken <-> Filecoin (row 32 now holds Filecoin, row 33 now holds RenderToken)
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.23"
$ws.Range("E32").Value = "  +3.68%  "

$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").Value = "'4.93"
$ws.Range("E33").Value = "  +5.51%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  -0.01%  "

# Row 35 - Celestia
$ws.Range("D35").Value = "'18.41"

# Row 36 - WEMIXToken
$ws.Range("E36").Value = "  -0.19%  "

# Row 37 - Hedera
$ws.Range("E37").Value = "  -0.01%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  +0.27%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  +0.61%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +1.30%  "

# Row 41 - LidoDAOToken
$ws.Range("E41").Value = "  -2.35%  "

# Row 42 - Maker
$ws.Range("D42").Value = "1.997.31"
$ws.Range("E42").Value = "  -0.21%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  -0.26%  "

# Row 44 - ApeXProtocol
$ws.Range("E44").Value = "  -5.74%  "

# Row 45 - FraxShare
$ws.Range("E45").Value = "  -1.82%  "

# Row 46 - EnergySwap
$ws.Range("D46").Value = "'17.48"
$ws.Range("E46").Value = "  -0.85%  "

# Row 47 - NEARProtocol
$ws.Range("D47").Value = "'2.84"
$ws.Range("E47").Value = "  -1.21%  "

# Row 48 - MultiversX
$ws.Range("D48").Value = "'54.74"
$ws.Range("E48").Value = "  -2.55%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "2.539.36"
$ws.Range("E49").Value = "  +0.79%  "

# Row 50 - BitcoinSV
$ws.Range("D50").Value = "'73.94"
$ws.Range("E50").Value = "  +5.06%  "

# Row 51 - Stacks
$ws.Range("E51").Value = "  +0.95%  "
